$wb = $excel.ActiveWorkbook

# ===== Overall =====
$ws = $wb.Worksheets.Item('Overall')
$ws.Cells.Item(1, 1).Value = 'Share of 990 filers with government grants at risk'
$headers = @('Number of 990 filers with government grants', 'Total government grants ($)', 'Size of operating surplus with government grants', 'Size of operating surplus without government grants')
for ($c = 0; $c -lt $headers.Length; $c++) { $ws.Cells.Item(1, $c + 2).Value = $headers[$c] }
$rows = @(
    ,@('''69.47%', '''2,532', '''$5,748,555,291', '''8.71%', '''-15.65%')
)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ===== County =====
$ws = $wb.Worksheets.Item('County')
$ws.Cells.Item(1, 1).Value = 'Geography'
$headers = @('Share of 990 filers with government grants at risk', 'Number of 990 filers with government grants', 'Total government grants ($)', 'Size of operating surplus with government grants', 'Size of operating surplus without government grants')
for ($c = 0; $c -lt $headers.Length; $c++) { $ws.Cells.Item(1, $c + 2).Value = $headers[$c] }
$rows = @(
    @('United States', '''67.35%', '''103,475', '''$267,700,640,005', '''9.05%', '''-12.83%'),
    @('New Jersey', '''69.47%', '''2,532', '''$5,748,555,291', '''8.71%', '''-15.65%'),
    @('Atlantic County', '''69.70%', '''66', '''$68,347,256', '''12.28%', '''-16.18%'),
    @('Bergen County', '''61.90%', '''231', '''$421,652,415', '''8.45%', '''-9.28%'),
    @('Burlington County', '''60.22%', '''93', '''$169,217,538', '''8.73%', '''-7.07%'),
    @('Camden County', '''75.35%', '''142', '''$356,399,434', '''5.59%', '''-28.24%'),
    @('Cape May County', '''63.16%', '''38', '''$18,248,540', '''18.77%', '''-7.23%'),
    @('Cumberland County', '''75.00%', '''44', '''$131,323,024', '''7.62%', '''-34.15%'),
    @('Essex County', '''73.93%', '''280', '''$1,090,505,242', '''8.36%', '''-25.48%'),
    @('Gloucester County', '''67.57%', '''37', '''$74,078,240', '''6.03%', '''-18.35%'),
    @('Hudson County', '''78.81%', '''118', '''$354,570,387', '''8.68%', '''-33.48%'),
    @('Hunterdon County', '''73.02%', '''63', '''$29,156,955', '''8.95%', '''-22.08%'),
    @('Mercer County', '''68.54%', '''213', '''$706,728,793', '''9.41%', '''-10.28%'),
    @('Middlesex County', '''77.27%', '''176', '''$607,737,194', '''4.69%', '''-23.97%'),
    @('Monmouth County', '''67.98%', '''203', '''$578,637,752', '''9.19%', '''-11.47%'),
    @('Morris County', '''65.48%', '''168', '''$182,317,810', '''12.29%', '''-10.68%'),
    @('Ocean County', '''62.30%', '''183', '''$177,655,803', '''9.65%', '''-8.85%'),
    @('Passaic County', '''70.16%', '''124', '''$289,486,834', '''5.89%', '''-18.99%'),
    @('Salem County', '''82.61%', '''23', '''$23,512,824', '''4.84%', '''-36.71%'),
    @('Somerset County', '''70.00%', '''120', '''$141,224,114', '''9.72%', '''-11.85%'),
    @('Sussex County', '''63.89%', '''36', '''$16,935,365', '''12.05%', '''-8.81%'),
    @('Union County', '''69.40%', '''134', '''$265,256,317', '''9.78%', '''-15.85%'),
    @('Warren County', '''72.50%', '''40', '''$45,563,454', '''8.40%', '''-19.94%')
)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ===== Congressional District =====
$ws = $wb.Worksheets.Item('Congressional District')
$ws.Cells.Item(1, 1).Value = 'Geography'
$headers = @('Share of 990 filers with government grants at risk', 'Number of 990 filers with government grants', 'Total government grants ($)', 'Size of operating surplus with government grants', 'Size of operating surplus without government grants')
for ($c = 0; $c -lt $headers.Length; $c++) { $ws.Cells.Item(1, $c + 2).Value = $headers[$c] }
$rows = @(
    @('United States', '''67.35%', '''103,475', '''$267,700,640,005', '''9.05%', '''-12.83%'),
    @('New Jersey', '''69.47%', '''2,532', '''$5,748,555,291', '''8.71%', '''-15.65%'),
    @('Congressional District 1', '''74.42%', '''172', '''$428,033,965', '''5.81%', '''-27.36%'),
    @('Congressional District 10', '''78.79%', '''231', '''$1,115,819,657', '''6.85%', '''-36.69%'),
    @('Congressional District 11', '''62.15%', '''251', '''$339,142,904', '''12.34%', '''-7.66%'),
    @('Congressional District 12', '''69.82%', '''275', '''$794,893,931', '''8.73%', '''-11.78%'),
    @('Congressional District 2', '''71.64%', '''201', '''$266,696,450', '''9.25%', '''-16.33%'),
    @('Congressional District 3', '''65.31%', '''196', '''$353,397,130', '''9.15%', '''-10.77%'),
    @('Congressional District 4', '''61.21%', '''232', '''$244,801,784', '''10.08%', '''-8.10%'),
    @('Congressional District 5', '''63.60%', '''228', '''$414,819,133', '''9.09%', '''-9.97%'),
    @('Congressional District 6', '''75.24%', '''210', '''$925,702,846', '''6.88%', '''-18.59%'),
    @('Congressional District 7', '''69.31%', '''277', '''$223,888,992', '''9.32%', '''-14.93%'),
    @('Congressional District 8', '''81.45%', '''124', '''$345,582,456', '''9.16%', '''-32.84%'),
    @('Congressional District 9', '''67.41%', '''135', '''$295,776,043', '''6.20%', '''-17.63%')
)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ===== Size =====
$ws = $wb.Worksheets.Item('Size')
$ws.Cells.Item(1, 1).Value = 'Size'
$headers = @('Share of 990 filers with government grants at risk', 'Number of 990 filers with government grants', 'Total government grants ($)', 'Size of operating surplus with government grants', 'Size of operating surplus without government grants')
for ($c = 0; $c -lt $headers.Length; $c++) { $ws.Cells.Item(1, $c + 2).Value = $headers[$c] }
$rows = @(
    @('Between $100K and $499K', '''68.34%', '''679', '''$72,988,085', '''11.19%', '''-17.42%'),
    @('Between $1M and $4.99M', '''72.09%', '''713', '''$705,050,939', '''8.22%', '''-17.63%'),
    @('Between $500K and $999K', '''67.73%', '''375', '''$114,147,311', '''8.02%', '''-16.02%'),
    @('Between $5M and $9.99M', '''70.20%', '''245', '''$741,154,506', '''8.65%', '''-16.00%'),
    @('Greater than $10M', '''68.64%', '''338', '''$4,107,627,504', '''5.94%', '''-8.43%'),
    @('Less than $100K', '''67.58%', '''182', '''$7,586,946', '''14.95%', '''-16.57%'),
    @('Total', '''69.47%', '''2,532', '''$5,748,555,291', '''8.71%', '''-15.65%')
)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ===== Subsector =====
$ws = $wb.Worksheets.Item('Subsector')
$ws.Cells.Item(1, 1).Value = 'Subsector'
$headers = @('Share of 990 filers with government grants at risk', 'Number of 990 filers with government grants', 'Total government grants ($)', 'Size of operating surplus with government grants', 'Size of operating surplus without government grants')
for ($c = 0; $c -lt $headers.Length; $c++) { $ws.Cells.Item(1, $c + 2).Value = $headers[$c] }
$rows = @(
    @('Arts, Culture, and Humanities', '''71.07%', '''197', '''$131,126,992', '''14.89%', '''-13.82%'),
    @('Education (Excluding Universities)', '''70.37%', '''324', '''$1,093,043,901', '''8.71%', '''-15.28%'),
    @('Environment and Animals', '''50.68%', '''73', '''$23,114,679', '''18.46%', '''-0.06%'),
    @('Health (Excluding Hospitals)', '''72.88%', '''306', '''$534,237,726', '''6.75%', '''-23.37%'),
    @('Hospitals', '''68.42%', '''19', '''$730,556,540', '''3.78%', '''-3.90%'),
    @('Human Services', '''70.08%', '''839', '''$1,107,531,305', '''8.37%', '''-18.72%'),
    @('International, Foreign Affairs', '''35.29%', '''17', '''$2,759,441', '''15.80%', '''7.52%'),
    @('Mutual/Membership Benefit', '''0.00%', '''1', '''$5,900', '''44.52%', '''41.92%'),
    @('Public, Societal Benefit', '''73.65%', '''167', '''$190,428,984', '''7.54%', '''-28.66%'),
    @('Religion Related', '''50.00%', '''64', '''$36,677,987', '''12.88%', '''-0.13%'),
    @('Unclassified', '''69.92%', '''502', '''$1,215,602,157', '''8.68%', '''-16.34%'),
    @('Universities', '''78.26%', '''23', '''$683,469,679', '''3.04%', '''-9.32%'),
    @('Total', '''69.47%', '''2,532', '''$5,748,555,291', '''8.71%', '''-15.65%')
)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

